$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.837.57"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.615.16"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.24"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.74"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "1.840.30"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "1.615.38"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").Value = "26.853.08"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.84"
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.50"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("E23").Value = "  -7.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.47"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.48"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -4.24%  "
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.716"
$ws.Range("E33").Value = "  +30.94%  "
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").Value = "1.323.69"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.829"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.47"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").Value = "1.753.79"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.31"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.820"
$ws.Range("E48").Value = "  +8.10%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0981"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0970"
$ws.Range("E51").Value = "  -7.91%  "
